$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data starts at row 2, header at row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C holds a "Förändrad" (changed) date for every data row; update it
# from 2023-09-10 (serial 45179) to 2023-09-11 (serial 45180) for all rows.
$ws.Range("C2:C$lastRow").Value = 45180
